$wb = $excel.ActiveWorkbook

# --- Update the conversion summary text on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.94 = 32503.97 pesos`n✅ 32503.97 pesos = 7.89 = 953.19 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 126
$ws2.Range("O10").Value = 4095.5
$ws2.Range("N12").Value = 4120
$ws2.Range("O12").Value = 120.82
